$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 86-87; everything currently at row 86 onward
# shifts down by two (old 86 -> new 88, ..., old 97 -> new 99).
$ws.Range("A86:A87").EntireRow.Insert()

# New row 86
$ws.Range("A86").Value = 4
$ws.Range("B86").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C86").Value = "Los Lagos"
$ws.Range("D86").Value = 44476
$ws.Range("E86").Value = 10
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100101
$ws.Range("H86").Value = "Berries"
$ws.Range("I86").Value = 100112025
$ws.Range("J86").Value = "Frutilla"
$ws.Range("K86").Value = "Sin especificar"
$ws.Range("L86").Value = "Especial"
$ws.Range("M86").Value = 200
$ws.Range("N86").Value = 15000
$ws.Range("O86").Value = 15000
$ws.Range("P86").Value = 15000
$ws.Range("Q86").Value = "$/bandeja 7 kilos"
$ws.Range("R86").Value = "Provincia de Melipilla"
$ws.Range("S86").Value = 2143
$ws.Range("T86").Value = 7

# New row 87
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = 44476
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100101
$ws.Range("H87").Value = "Berries"
$ws.Range("I87").Value = 100112025
$ws.Range("J87").Value = "Frutilla"
$ws.Range("K87").Value = "Sin especificar"
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 200
$ws.Range("N87").Value = 13000
$ws.Range("O87").Value = 13000
$ws.Range("P87").Value = 13000
$ws.Range("Q87").Value = "$/bandeja 7 kilos"
$ws.Range("R87").Value = "Provincia de Melipilla"
$ws.Range("S87").Value = 1857
$ws.Range("T87").Value = 7
